$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-7 from 2023-09-05 (45174) to 2023-09-06 (45175)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45175
}
